$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D (shifts old D -> E, old E -> F);
# the inserted column inherits the header style from column D.
$ws.Columns.Item(4).Insert()

# New header for the inserted column
$ws.Cells.Item(1, 4).Value = "other_hypercap_threshold"

# Fill in the new column's values (other_hypercap_threshold)
$ws.Cells.Item(2, 4).Value = 8.4
$ws.Cells.Item(3, 4).Value = 7.56
$ws.Cells.Item(4, 4).Value = 8.029999999999999
$ws.Cells.Item(5, 4).Value = 10.33
$ws.Cells.Item(6, 4).Value = 13.49
$ws.Cells.Item(7, 4).Value = 2.32
$ws.Cells.Item(8, 4).Value = 6.2
$ws.Cells.Item(9, 4).Value = 6.27
$ws.Cells.Item(10, 4).Value = 16.15
$ws.Cells.Item(11, 4).Value = 18.84
$ws.Cells.Item(12, 4).Value = 2.41

# Corrected abg_hypercap_threshold value (column B)
$ws.Cells.Item(4, 2).Value = 8.779999999999999

# Corrected pco2_threshold_any values (column E, formerly D)
$ws.Cells.Item(4, 5).Value = 7.81
$ws.Cells.Item(5, 5).Value = 10.12
$ws.Cells.Item(6, 5).Value = 14.03
$ws.Cells.Item(10, 5).Value = 15.97
$ws.Cells.Item(11, 5).Value = 19.46

# Corrected vbg_hypercap_threshold values (column F, formerly E)
$ws.Cells.Item(4, 6).Value = 6.44
$ws.Cells.Item(6, 6).Value = 14.88
$ws.Cells.Item(9, 6).Value = 5.66
$ws.Cells.Item(10, 6).Value = 14.01
$ws.Cells.Item(11, 6).Value = 25.4
